# Update cryptocurrency price/volume data per the latest GitHub Actions scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''67.315.29'
$ws.Range('E2').Value = '  -2.18%  '
$ws.Range('D3').Value = '''3.240.03'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''583.19'
$ws.Range('E5').Value = '  -4.77%  '
$ws.Range('D6').Value = '''144.73'
$ws.Range('E6').Value = '  -12.79%  '
$ws.Range('D8').Value = '''3.217.64'
$ws.Range('E8').Value = '  -5.83%  '
$ws.Range('D9').Value = '''0.531'
$ws.Range('E9').Value = '  -10.58%  '
$ws.Range('D10').Value = '''0.165'
$ws.Range('E10').Value = '  -14.61%  '
$ws.Range('D11').Value = '''6.72'
$ws.Range('E11').Value = '  -2.67%  '
$ws.Range('D12').Value = '''0.491'
$ws.Range('E12').Value = '  -12.14%  '
$ws.Range('D13').Value = '''0.0000240'
$ws.Range('E13').Value = '  -10.32%  '
$ws.Range('D14').Value = '''37.02'
$ws.Range('E14').Value = '  -15.36%  '
$ws.Range('D15').Value = '''3.760.31'
$ws.Range('E15').Value = '  -5.45%  '
$ws.Range('D16').Value = '''67.368.32'
$ws.Range('E16').Value = '  -2.20%  '
$ws.Range('D17').Value = '''3.245.57'
$ws.Range('E17').Value = '  -5.24%  '
$ws.Range('E18').Value = '  -6.47%  '
$ws.Range('D19').Value = '''6.92'
$ws.Range('E19').Value = '  -14.23%  '
$ws.Range('D20').Value = '''508.17'
$ws.Range('E20').Value = '  -11.58%  '
$ws.Range('D21').Value = '''14.54'
$ws.Range('E21').Value = '  -14.29%  '
$ws.Range('D22').Value = '''0.735'
$ws.Range('E22').Value = '  -12.41%  '
$ws.Range('D23').Value = '''7.54'
$ws.Range('E23').Value = '  -15.71%  '
$ws.Range('D24').Value = '''83.80'
$ws.Range('D25').Value = '''13.04'
$ws.Range('E25').Value = '  -12.74%  '
$ws.Range('D26').Value = '''0.999'
$ws.Range('E26').Value = '  -0.14%  '
$ws.Range('D27').Value = '''3.15'
$ws.Range('E27').Value = '  -12.73%  '
$ws.Range('D28').Value = '''2.08'
$ws.Range('E28').Value = '  -12.79%  '
$ws.Range('D29').Value = '''7.74'
$ws.Range('E29').Value = '  -8.74%  '
$ws.Range('D30').Value = '''28.23'
$ws.Range('E30').Value = '  -13.00%  '
$ws.Range('E31').Value = '  -5.26%  '
$ws.Range('E32').Value = '  -6.33%  '
$ws.Range('D33').Value = '''6.32'
$ws.Range('E33').Value = '  -18.62%  '
$ws.Range('E34').Value = '  +0.00%  '
$ws.Range('D35').Value = '''5.51'
$ws.Range('E35').Value = '  -15.29%  '
$ws.Range('D36').Value = '''55.08'
$ws.Range('E36').Value = '  -1.81%  '
$ws.Range('D37').Value = '''498.32'
$ws.Range('E37').Value = '  -14.40%  '
$ws.Range('E38').Value = '  -7.83%  '
$ws.Range('D39').Value = '''0.0828'
$ws.Range('E39').Value = '  -12.54%  '
$ws.Range('B40').Value = 'Cosmos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D40').Value = '''8.69'
$ws.Range('E40').Value = '  -16.44%  '
$ws.Range('B41').Value = 'Kaspa'
$ws.Range('C41').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D41').Value = '''0.122'
$ws.Range('E41').Value = '  -12.55%  '
$ws.Range('D42').Value = '''2.883.49'
$ws.Range('E42').Value = '  -10.10%  '
$ws.Range('D43').Value = '''2.64'
$ws.Range('E43').Value = '  -13.54%  '
$ws.Range('D44').Value = '''0.258'
$ws.Range('E44').Value = '  -11.23%  '
$ws.Range('E45').Value = '  -0.10%  '
$ws.Range('D46').Value = '''2.12'
$ws.Range('E46').Value = '  -10.05%  '
$ws.Range('D47').Value = '''25.83'
$ws.Range('E47').Value = '  -16.67%  '
$ws.Range('D48').Value = '''0.0₃0544'
$ws.Range('E48').Value = '  -18.84%  '
$ws.Range('D49').Value = '''122.97'
$ws.Range('E49').Value = '  -7.00%  '
$ws.Range('E50').Value = '  -11.55%  '
$ws.Range('D51').Value = '''2.23'
$ws.Range('E51').Value = '  -19.19%  '
